$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update test data values on the Smoke_Suite sheet ---
# Row 2 (TestCreateChangeTicket): Change ID -> CHG0037654, Status -> Passed
$ws.Range("C2").Value = "CHG0037654"
$ws.Range("E2").Value = "Passed"

# Row 3 (TestUpdateChangeTicket): Change ID -> Passed (per source data)
# the cell previously had no value (just a bordered/empty style) -- clear the
# formatting before writing so it matches the plain unstyled string cell
$ws.Range("C3").Clear()
$ws.Range("C3").Value = "Passed"

# Row 4 (TestApproveChangeTicket): Change ID -> CHG0037654, Status -> Passed
$ws.Range("C4").Clear()
$ws.Range("C4").Value = "CHG0037654"
$ws.Range("E4").Clear()
$ws.Range("E4").Value = "Passed"

# --- Update view/selection state ---
[void]$ws.Range("C2").Select()
